# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last row previously had a "date only" number format to visually mark
# it as the latest entry; once a new row is appended, that formatting moves
# to the new last row, and the old last row reverts to the standard
# date+time format used by all other data rows.
$ws.Cells.Item(84, 1).NumberFormat = $ws.Cells.Item(83, 1).NumberFormat

# Append the new day's data as row 85
$ws.Cells.Item(85, 1).Value = 45825
$ws.Cells.Item(85, 2).Value = 360
$ws.Cells.Item(85, 3).Value = 365
$ws.Cells.Item(85, 4).Value = 367

# Mark the newly appended row as the latest entry with the date-only format
$ws.Cells.Item(85, 1).NumberFormat = "YYYY-MM-DD"
